$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.378.41"
$ws.Range("E2").Value = "  +2.40%  "
$ws.Range("D3").Value = "'2.035.71"
$ws.Range("E3").Value = "  +4.35%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'248.64"
$ws.Range("E5").Value = "  +2.38%  "
$ws.Range("D6").Value = "'0.630"
$ws.Range("E6").Value = "  +2.42%  "
$ws.Range("D7").Value = "'60.96"
$ws.Range("E7").Value = "  +1.61%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +6.12%  "
$ws.Range("E10").Value = "  +3.71%  "
$ws.Range("E11").Value = "  +2.28%  "
$ws.Range("D12").Value = "'15.39"
$ws.Range("E12").Value = "  +8.81%  "
$ws.Range("D13").Value = "'0.864"
$ws.Range("E13").Value = "  +4.45%  "
$ws.Range("D14").Value = "'2.336.83"
$ws.Range("E14").Value = "  +4.29%  "
$ws.Range("D15").Value = "'22.52"
$ws.Range("E15").Value = "  +4.89%  "
$ws.Range("D16").Value = "'5.55"
$ws.Range("E16").Value = "  +6.36%  "
$ws.Range("D17").Value = "'2.033.45"
$ws.Range("E17").Value = "  +4.39%  "
$ws.Range("D18").Value = "'37.318.37"
$ws.Range("E18").Value = "  +2.56%  "
$ws.Range("D19").Value = "'70.86"
$ws.Range("E19").Value = "  +2.25%  "
$ws.Range("D20").Value = "'0.0₃0871"
$ws.Range("E20").Value = "  +2.87%  "
$ws.Range("E21").Value = "  +4.55%  "
$ws.Range("D22").Value = "'231.60"
$ws.Range("E22").Value = "  +1.32%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'2.57"
$ws.Range("E24").Value = "  +5.85%  "
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("D26").Value = "'9.53"
$ws.Range("E26").Value = "  +4.30%  "
$ws.Range("D27").Value = "'164.00"
$ws.Range("E27").Value = "  +2.71%  "
$ws.Range("E28").Value = "  -3.77%  "
$ws.Range("D29").Value = "'19.88"
$ws.Range("E29").Value = "  +3.54%  "
$ws.Range("E30").Value = "  +5.95%  "
$ws.Range("E31").Value = "  +3.06%  "
$ws.Range("E32").Value = "  +2.92%  "
$ws.Range("D33").Value = "'0.0673"
$ws.Range("E33").Value = "  +10.45%  "
$ws.Range("E34").Value = "  +2.89%  "
$ws.Range("D35").Value = "'2.52"
$ws.Range("E35").Value = "  +11.48%  "
$ws.Range("D36").Value = "'3.68"
$ws.Range("E36").Value = "  +6.96%  "
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("D39").Value = "'5.48"
$ws.Range("E39").Value = "  +0.95%  "
$ws.Range("B40").Value = "HuobiToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D40").Value = "'3.00"
$ws.Range("E40").Value = "  +3.11%  "
$ws.Range("B41").Value = "Cronos"
$ws.Range("C41").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D41").Value = "'0.0985"
$ws.Range("E41").Value = "  +3.45%  "
$ws.Range("D42").Value = "'17.20"
$ws.Range("E42").Value = "  +9.40%  "
$ws.Range("E43").Value = "  +2.99%  "
$ws.Range("D44").Value = "'0.0215"
$ws.Range("E44").Value = "  +3.36%  "
$ws.Range("D45").Value = "'93.18"
$ws.Range("E45").Value = "  +5.41%  "
$ws.Range("E46").Value = "  +4.59%  "
$ws.Range("D47").Value = "'1.391.06"
$ws.Range("E47").Value = "  +2.53%  "
$ws.Range("D48").Value = "'7.52"
$ws.Range("E48").Value = "  +6.06%  "
$ws.Range("D49").Value = "'2.20"
$ws.Range("E49").Value = "  +21.35%  "
$ws.Range("D50").Value = "'2.86"
$ws.Range("E50").Value = "  +1.36%  "
$ws.Range("D51").Value = "'46.48"
$ws.Range("E51").Value = "  +2.27%  "
